$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "triple" block: BNB/BTC, SC/BTC, SC/BNB labels + values + derived formulas
$ws.Range("A22").Value = "BNB/BTC"
$ws.Range("C22").Value = "SC/BTC"
$ws.Range("E22").Value = "SC/BNB"

$ws.Range("A24").Value = 0.0045269999999999998
$ws.Range("C24").Value = 0.00000035999999999999999
$ws.Range("E24").Value = 0.000079099999999999998

$ws.Range("A26").Formula = "=A24"
$ws.Range("C26").Formula = "=A26/C24"
$ws.Range("E26").Formula = "=C26*E24"

# Update the active selection to match the author's final cursor position
$ws.Range("F14").Select()
